$wb = $excel.ActiveWorkbook

# --- Add Athlete sheet after the last existing sheet (Team) ---
$athlete = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$athlete.Name = "Athlete"

$athleteData = @(
    @("Hermione Granger", 1, 4, 7, "Gryfindor", 0, 1, 0),
    @("Ron Weasley", 2, 4, 8, "Hufflepuff", 1, 1, 1),
    @("Albus Dumbledore", 4, 6, 9, "Ravenclaw", 2, 2, 0),
    @("Ginny Weasley", 8, 6, 1, "Slytherin", 0, 3, 0),
    @("Rubeus Hagrid", 16, 5, 12, "Gryfindor", 2, 3, 0),
    @("Minevra McGonagall", 32, 6, 5, "Hufflepuff", 1, 2, 0)
)

for ($r = 0; $r -lt $athleteData.Length; $r++) {
    $row = $athleteData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $athlete.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Names in column A are long - widen it to fit (matches bestFit column behavior)
$athlete.Columns.Item(1).AutoFit() | Out-Null

# --- Add Rides sheet after Athlete ---
$rides = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$rides.Name = "Rides"

$ridesData = @(
    @(1, 1),
    @(12, 1)
)

for ($r = 0; $r -lt $ridesData.Length; $r++) {
    $row = $ridesData[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $rides.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# --- Set selections per sheet (match target view state) ---
$wb.Worksheets.Item("Team").Select() | Out-Null
$wb.Worksheets.Item("Team").Range("C2").Select() | Out-Null

$athlete.Select() | Out-Null
$athlete.Range("H7").Select() | Out-Null

$rides.Select() | Out-Null
$rides.Range("B3").Select() | Out-Null

